$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: B2 gets the new JIRA id, C2 gets the new description.
# A2 (TCID) and D2 (Runmode) remain unchanged.
$ws.Range("B2").Value = "WAT-231"
$ws.Range("C2").Value = "Verify that Sign in using email and password (Steam Login)has been added to that login page"

# Update the active selection to C2 (was A2).
$ws.Range("C2").Select()
